$wb = $excel.ActiveWorkbook

# start_price sheet
$ws = $wb.Worksheets.Item("start_price")
$ws.Range("A2").Value = 45.81

# Linear sheet
$ws = $wb.Worksheets.Item("Linear")
$ws.Range("B2").Value = 0.002615313252274685
$ws.Range("B3").Value = -0.07923649046753914
$ws.Range("B4").Value = 1.808812253881741
$ws.Range("B5").Value = "[1.0, 0.2254609248719924, 0.011905922511050411, -0.04650682427590289, -0.044098553131532, 0.03033564586727745, 0.20514982032503615, 0.3877108396390888, 0.1942447139389808, 0.009528078106427687, -0.04595394384375145, -0.056481618338936675, 0.01485539905854682, 0.19434317132285397, 0.361201072961821, 0.18215981285519614, 0.0002690086851011737, -0.041353053919284066, -0.0427252927041747, -0.007041783067033135]"

# NonLinear sheet
$ws = $wb.Worksheets.Item("NonLinear")
$ws.Range("B3").Value = 0.9084993899959333
$ws.Range("B4").Value = 0.01240513083201955
$ws.Range("B5").Value = -0.02343746870349536
$ws.Range("B6").Value = 1.936312310869361
$ws.Range("B7").Value = 0.08591928392187169
$ws.Range("B8").Value = -0.250423008170013
$ws.Range("B9").Value = 1.690068338648737
$ws.Range("B10").Value = "[1.0, 0.2230198756841849, 0.012011688183475427, -0.045039170332152725, -0.04377016895457832, 0.030704865541344748, 0.20223488596585942, 0.3870270433771894, 0.19101749491604164, 0.008847830227062035, -0.04559331156822698, -0.05591871818472744, 0.013269880549708563, 0.1923759249067009, 0.35924813719237947, 0.18161450219297498, -0.00011559649212930185, -0.04198627101776983, -0.04050061082577633, -0.007781171757771893]"
